# feat: add 2022-Q1 data
#
# The old "总计" (totals) sheet is repurposed in-place to become the new
# "2022-Q1" holdings sheet (keeps its original sheetId=3), and a brand new
# "总计" sheet is appended right after it (sheetId=4), rebuilt with the old
# totals rows plus a new leading row for the 2022-Q1 quarter.

$wb = $excel.ActiveWorkbook

$q3 = $wb.Worksheets.Item(2)   # "2021-Q3" - used as a formatting template
$q1 = $wb.Worksheets.Item(3)   # currently "总计", becomes "2022-Q1"

# ---------------------------------------------------------------------
# 1) Turn the existing "总计" sheet into the new "2022-Q1" holdings sheet.
# ---------------------------------------------------------------------
$q1.Cells.Clear()
$q1.Name = "2022-Q1"

# Borrow the header-row / index-column formatting (bold, bordered,
# centered) from the "2021-Q3" sheet, which already uses it, instead of
# creating new style entries.
$q3.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$q3.Range("A2:A4").Copy()
$q1.Range("A2:A4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Header row (row 1)
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row 2
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "'000593"
$q1.Range("C2").Value = "易方达标普全球高端消费品指数增强(QDII)-美元现汇"
$q1.Range("D2").Value = "'1.93"
$q1.Range("E2").Value = "'92.46"
$q1.Range("F2").Value = "'8.46"
$q1.Range("G2").Value = "'0.1633"
$q1.Range("H2").Value = 3

# Row 3
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "'005676"
$q1.Range("C3").Value = "易方达标普全球高端消费品指数增强C(QDII) - 人民币"
$q1.Range("D3").Value = "'1.93"
$q1.Range("E3").Value = "'92.46"
$q1.Range("F3").Value = "'8.46"
$q1.Range("G3").Value = "'0.1633"
$q1.Range("H3").Value = 3

# Row 4
$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "'118002"
$q1.Range("C4").Value = "易方达标普全球高端消费品指数增强A(QDII) - 人民币"
$q1.Range("D4").Value = "'1.93"
$q1.Range("E4").Value = "'92.46"
$q1.Range("F4").Value = "'8.46"
$q1.Range("G4").Value = "'0.1633"
$q1.Range("H4").Value = 3

# ---------------------------------------------------------------------
# 2) Add a fresh "总计" sheet right after "2022-Q1", rebuilt with the old
#    totals rows plus a new leading row for 2022-Q1.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

# Match the outline settings (summary row below / summary column on the
# right) used by the sibling sheets.
$total.Outline.SummaryRow = 1
$total.Outline.SummaryColumn = 1

# Match the page margins used by the sibling sheets (0.75in/0.75in/1in/1in,
# 0.5in header/footer - points are 1/72 inch).
$total.PageSetup.LeftMargin = 54
$total.PageSetup.RightMargin = 54
$total.PageSetup.TopMargin = 72
$total.PageSetup.BottomMargin = 72
$total.PageSetup.HeaderMargin = 36
$total.PageSetup.FooterMargin = 36

# Borrow the same header/index-column styling again.
$q3.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$q3.Range("A2:A4").Copy()
$total.Range("A2:A4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Header row
$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

# Row 2 - new 2022-Q1 entry
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.49

# Row 3 - previously row 2 (2021-Q3)
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q3"
$total.Range("C3").Value = 6
$total.Range("D3").Value = 4.49

# Row 4 - previously row 3 (2021-Q2)
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q2"
$total.Range("C4").Value = 8
$total.Range("D4").Value = 6.71
